$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "40V (A)" header label (used by H8, H9 and H10 which share the same
# string) needs to become "40V Rail(A)" per the new implementation.
$ws.Range("H8").Value = "40V Rail(A)"
$ws.Range("H9").Value = "40V Rail(A)"
$ws.Range("H10").Value = "40V Rail(A)"

# Update the sheet selection/view: select H9:H10 (active cell H9) which also
# resets the scrolled top-left cell back to its default.
$ws.Range("H9:H10").Select()
